# Auto-generated cell updates applying the cryptos.xlsx diff (cryptocurrency price refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '20.183.10'
$ws.Cells.Item(2, 5).Value = '  +1.38%  '
$ws.Cells.Item(3, 4).Value = '1.442.94'
$ws.Cells.Item(3, 5).Value = '  +2.95%  '
$ws.Cells.Item(4, 4).Value = '''1.008'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.48%  '
$ws.Cells.Item(5, 4).Value = '''0.9123'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -9.00%  '
$ws.Cells.Item(6, 4).Value = '''277.35'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +2.81%  '
$ws.Cells.Item(7, 4).Value = '''0.3665'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +0.32%  '
$ws.Cells.Item(8, 4).Value = '''0.3132'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +3.46%  '
$ws.Cells.Item(9, 4).Value = '''38.97'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -0.32%  '
$ws.Cells.Item(10, 4).Value = '''1.020'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +5.48%  '
$ws.Cells.Item(11, 4).Value = '''0.06526'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +1.38%  '
$ws.Cells.Item(12, 4).Value = '''1.002'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -0.09%  '
$ws.Cells.Item(13, 4).Value = '''5.395'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +2.76%  '
$ws.Cells.Item(14, 4).Value = '''17.53'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +6.03%  '
$ws.Cells.Item(15, 4).Value = '''6.077'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +0.48%  '
$ws.Cells.Item(16, 2).Value = 'ShibaInu'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(16, 4).Value = '''0.00001015'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +1.82%  '
$ws.Cells.Item(17, 2).Value = 'WrappedEther'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(17, 4).Value = '1.440.78'
$ws.Cells.Item(17, 5).Value = '  +2.53%  '
$ws.Cells.Item(18, 4).Value = '''0.9356'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -6.65%  '
$ws.Cells.Item(19, 4).Value = '''0.05636'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -0.45%  '
$ws.Cells.Item(20, 4).Value = '''67.78'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -4.83%  '
$ws.Cells.Item(21, 2).Value = 'Avalanche'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(21, 4).Value = '''14.47'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +2.61%  '
$ws.Cells.Item(22, 2).Value = 'Uniswap'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(22, 4).Value = '''5.401'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -1.28%  '
$ws.Cells.Item(23, 4).Value = '''10.81'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +2.58%  '
$ws.Cells.Item(24, 4).Value = '''2.257'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -0.77%  '
$ws.Cells.Item(25, 4).Value = '20.174.33'
$ws.Cells.Item(26, 4).Value = '''2.183'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -0.84%  '
$ws.Cells.Item(27, 4).Value = '''136.53'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +0.81%  '
$ws.Cells.Item(28, 5).Value = '  +2.75%  '
$ws.Cells.Item(29, 4).Value = '1.593.34'
$ws.Cells.Item(29, 5).Value = '  +1.81%  '
$ws.Cells.Item(30, 4).Value = '''110.35'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +2.89%  '
$ws.Cells.Item(31, 4).Value = '''3.785'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -1.08%  '
$ws.Cells.Item(32, 4).Value = '''0.8052'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +2.05%  '
$ws.Cells.Item(33, 4).Value = '''4.823'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -6.81%  '
$ws.Cells.Item(34, 4).Value = '''0.07704'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +1.83%  '
$ws.Cells.Item(35, 4).Value = '''0.05992'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +5.36%  '
$ws.Cells.Item(36, 4).Value = '''1.451'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +9.06%  '
$ws.Cells.Item(37, 4).Value = '''4.694'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +0.58%  '
$ws.Cells.Item(38, 4).Value = '''1.134'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +8.39%  '
$ws.Cells.Item(39, 4).Value = '''0.01998'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -0.20%  '
$ws.Cells.Item(40, 2).Value = 'Aptos'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(40, 4).Value = '''10.16'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +0.86%  '
$ws.Cells.Item(41, 2).Value = 'Frax'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(41, 4).Value = '''0.9322'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -6.89%  '
$ws.Cells.Item(42, 4).Value = '''0.1836'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -3.69%  '
$ws.Cells.Item(43, 4).Value = '''7.217'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -12.92%  '
$ws.Cells.Item(44, 4).Value = '''3.520'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +1.02%  '
$ws.Cells.Item(45, 4).Value = '''0.5234'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +0.71%  '
$ws.Cells.Item(46, 4).Value = '''12.04'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +1.00%  '
$ws.Cells.Item(47, 4).Value = '''118.94'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +9.33%  '
$ws.Cells.Item(48, 4).Value = '''0.5137'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +2.76%  '
$ws.Cells.Item(49, 4).Value = '''1.767'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +1.66%  '
$ws.Cells.Item(50, 4).Value = '''0.06318'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +3.25%  '
$ws.Cells.Item(51, 4).Value = '''0.9948'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -0.77%  '

Write-Output "Applied 110 cell updates"
